# Contact Tracing workbook update:
# - Cases table: refresh/replace sample rows with a larger (22-row) dataset and
#   make the table scrollable by expanding it to A1:J23.
# - Contacts table: drop the now-unused "Traced Date"/"Contacted date" values
#   for the first two data rows and refresh a couple of Id/CaseId/Added Date values.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# "Cases" sheet / table
# -----------------------------------------------------------------
$casesSheet = $wb.Worksheets.Item("Cases")

# Clear out the previous sample rows (2-5) completely; some old columns
# (H = "Traced Date") are not used by the new dataset at all.
$casesSheet.Range("A2:J5").ClearContents()

# Columns: Id, Test Date, Added Date, Postcode, Traced?, Dropped times,
#          Dropped?, Symptom date (column I - column H "Traced Date" unused)
$casesData = @(
    @(1, 44317, 44317.9528023727, "OX1", $false, 0, $false, "4/30/2021 12:00:00 AM"),
    @(2, 44317, 44317.9534171065, "OX1", $false, 0, $false, "5/1/2021 12:00:00 AM"),
    @(3, 44317, 44317.9625158681, "OX1", $false, 0, $false, $null),
    @(1002, 44327, 44327.7865993403, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1003, 44327, 44327.7869098611, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1004, 44327, 44327.7871325463, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1005, 44327, 44327.7874183102, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1006, 44327, 44327.7876739699, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1007, 44327, 44327.7878651042, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1008, 44327, 44327.7881653935, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1009, 44327, 44327.7885158333, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1010, 44327, 44327.7887260648, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1011, 44327, 44327.7889103241, "OX1", $false, 0, $false, "5/11/2021 12:00:00 AM"),
    @(1012, 44327, 44327.7896717014, "OX1", $false, 0, $false, $null),
    @(1013, 44327, 44327.7898875116, "OX1", $false, 0, $false, $null),
    @(1014, 44327, 44327.7900773148, "OX1", $false, 0, $false, $null),
    @(1015, 44327, 44327.7902346644, "OX1", $false, 0, $false, $null),
    @(1016, 44327, 44327.7905140393, "OX1", $false, 0, $false, $null),
    @(1017, 44327, 44327.7906821875, "OX1", $false, 0, $false, $null),
    @(1018, 44327, 44327.7908281944, "OX1", $false, 0, $false, $null),
    @(1019, 44327, 44327.7909767708, "OX1", $false, 0, $false, $null),
    @(1020, 44327, 44327.7911247338, "OX1", $false, 0, $false, $null)
)

$r = 2
foreach ($row in $casesData) {
    $casesSheet.Cells.Item($r, 1).Value = $row[0]
    $casesSheet.Cells.Item($r, 2).Value = $row[1]
    $casesSheet.Cells.Item($r, 3).Value = $row[2]
    $casesSheet.Cells.Item($r, 4).Value = $row[3]
    $casesSheet.Cells.Item($r, 5).Value = $row[4]
    $casesSheet.Cells.Item($r, 6).Value = $row[5]
    $casesSheet.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $casesSheet.Cells.Item($r, 9).Value = $row[7]
    }
    $r++
}

# Expand the "Cases" table/autofilter so the new rows are part of it and it
# becomes scrollable.
$casesTable = $casesSheet.ListObjects.Item("Cases")
$casesTable.Resize($casesSheet.Range("A1:J23"))

# -----------------------------------------------------------------
# "Contacts" sheet / table
# -----------------------------------------------------------------
$contactsSheet = $wb.Worksheets.Item("Contacts")

$contactsSheet.Cells.Item(2, 1).Value = 3
$contactsSheet.Cells.Item(2, 2).Value = 2
$contactsSheet.Cells.Item(2, 3).Value = 44319.7581985185
$contactsSheet.Cells.Item(2, 4).ClearContents()
$contactsSheet.Cells.Item(2, 5).ClearContents()

$contactsSheet.Cells.Item(3, 1).Value = 4
$contactsSheet.Cells.Item(3, 2).Value = 2
$contactsSheet.Cells.Item(3, 3).Value = 44319.7585063773
$contactsSheet.Cells.Item(3, 4).ClearContents()
$contactsSheet.Cells.Item(3, 5).ClearContents()

$wb.Save()
